# Student Role Accepting Script completed
# Applies the "ApplyRole / Feedback / UploadResume" onboarding update.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. RegistrationForm: swap demo student for the new one.
# ---------------------------------------------------------------------
$reg = $wb.Worksheets.Item("RegistrationForm")
$reg.Activate()
$reg.Range("A2").Value = "Cathalyin"
$reg.Range("B2").Value = "C"
$reg.Range("C2").Value = "priya.t+studentrolecreation3@icanio.com"
$reg.Range("C2").Select()

# ---------------------------------------------------------------------
# 2. EducationalDetails: no data change, just browsed through.
# ---------------------------------------------------------------------
$edu = $wb.Worksheets.Item("EducationalDetails")
$edu.Activate()
$edu.Range("U2").Select()

# ---------------------------------------------------------------------
# 3. AdditionalDetails: no data change, just browsed through.
# ---------------------------------------------------------------------
$additional = $wb.Worksheets.Item("AdditionalDetails")
$additional.Activate()
$additional.Range("F19").Select()

# ---------------------------------------------------------------------
# 4. PersonalDetails: gender + registration number update.
# ---------------------------------------------------------------------
$personal = $wb.Worksheets.Item("PersonalDetails")
$personal.Activate()
$personal.Range("A2").Value = 2345
$personal.Range("C2").Value = "Female"
$personal.Range("D7").Select()

# ---------------------------------------------------------------------
# 5. ProjectandInternship: skill4 update.
# ---------------------------------------------------------------------
$project = $wb.Worksheets.Item("ProjectandInternship")
$project.Activate()
$project.Range("N2").Value = "Selenium"
$project.Range("O2").Select()

# ---------------------------------------------------------------------
# 6. WorkExperience: no data change, just browsed through.
# ---------------------------------------------------------------------
$work = $wb.Worksheets.Item("WorkExperience")
$work.Activate()
$work.Range("U2").Select()

# ---------------------------------------------------------------------
# 7. CoursesAndTraining: skill1 update.
# ---------------------------------------------------------------------
$courses = $wb.Worksheets.Item("CoursesAndTraining")
$courses.Activate()
$courses.Range("H2").Value = "communication"
$courses.Range("C21").Select()

# ---------------------------------------------------------------------
# 8. ProofAndDocument: untouched, left as-is.
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 9. Add the three new sheets at the end of the workbook, in order:
#    ApplyRole, Feedback, UploadResume.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$applyRole = $wb.Worksheets.Add($null, $lastSheet)
$applyRole.Name = "ApplyRole"
$applyRole.Range("A1").Value = "search"
$applyRole.Range("O1").Value = "search"
$applyRole.Range("A2").Value = "Quality Roles"
$applyRole.Range("O2").Value = "Database Role"

$feedback = $wb.Worksheets.Add($null, $applyRole)
$feedback.Name = "Feedback"
$feedback.Range("A1").Value = "reference"
$feedback.Range("B1").Value = "feedback"
$feedback.Range("A2").Value = "PluginLive"
$feedback.Range("B2").Value = "Good and User Friendly"

$uploadResume = $wb.Worksheets.Add($null, $feedback)
$uploadResume.Name = "UploadResume"
$uploadResume.Range("A1").Value = "resumefile"
$uploadResume.Range("A2").Value = "C:\Users\ICANIO-10090\Pictures\Resume.jpg"

# ---------------------------------------------------------------------
# 10. Finish with ApplyRole as the active/selected sheet (activeTab=10).
# ---------------------------------------------------------------------
$applyRole.Activate()
$applyRole.Range("O2").Select()
